$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '''27.144.27'
$ws.Range('E2').Value = '''  +0.27%  '

# Row 3
$ws.Range('D3').Value = '''1.820.82'
$ws.Range('E3').Value = '''  -0.30%  '

# Row 4
$ws.Range('D4').Value = '''0.9993'
$ws.Range('E4').Value = '''  -0.20%  '

# Row 5
$ws.Range('D5').Value = '''312.40'
$ws.Range('E5').Value = '''  +0.26%  '

# Row 6
$ws.Range('D6').Value = '''0.9999'
$ws.Range('E6').Value = '''  -0.11%  '

# Row 7
$ws.Range('D7').Value = '''0.4466'
$ws.Range('E7').Value = '''  +5.34%  '

# Row 8
$ws.Range('D8').Value = '''0.3744'
$ws.Range('E8').Value = '''  +1.95%  '

# Row 9
$ws.Range('D9').Value = '''0.07478'
$ws.Range('E9').Value = '''  +3.42%  '

# Row 10
$ws.Range('D10').Value = '''0.8727'
$ws.Range('E10').Value = '''  +3.13%  '

# Row 11
$ws.Range('D11').Value = '''20.91'
$ws.Range('E11').Value = '''  +1.05%  '

# Row 12
$ws.Range('D12').Value = '''1.807.16'
$ws.Range('E12').Value = '''  -1.05%  '

# Row 13
$ws.Range('D13').Value = '''6.739'
$ws.Range('E13').Value = '''  +1.13%  '

# Row 14
$ws.Range('D14').Value = '''94.52'
$ws.Range('E14').Value = '''  +5.44%  '

# Row 15
$ws.Range('D15').Value = '''5.360'
$ws.Range('E15').Value = '''  +1.23%  '

# Row 16
$ws.Range('D16').Value = '''0.07105'
$ws.Range('E16').Value = '''  +0.92%  '

# Row 17
$ws.Range('D17').Value = '''0.9998'
$ws.Range('E17').Value = '''  -0.28%  '

# Row 18
$ws.Range('D18').Value = '''0.000008763'
$ws.Range('E18').Value = '''  +0.00%  '

# Row 19
$ws.Range('D19').Value = '''0.9996'
$ws.Range('E19').Value = '''  -0.15%  '

# Row 20
$ws.Range('D20').Value = '''15.02'
$ws.Range('E20').Value = '''  +1.02%  '

# Row 21
$ws.Range('D21').Value = '''27.162.91'
$ws.Range('E21').Value = '''  +0.25%  '

# Row 22
$ws.Range('D22').Value = '''5.233'
$ws.Range('E22').Value = '''  +1.90%  '

# Row 23
$ws.Range('D23').Value = '''10.95'
$ws.Range('E23').Value = '''  +1.31%  '

# Row 24
$ws.Range('B24').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C24').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D24').Value = '''2.041.71'
$ws.Range('E24').Value = '''  -0.45%  '

# Row 25
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Value = '''1.990'
$ws.Range('E25').Value = '''  +0.67%  '

# Row 26
$ws.Range('B26').Value = 'LidoDAOToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D26').Value = '''2.431'
$ws.Range('E26').Value = '''  +7.99%  '

# Row 27
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').Value = '''151.67'
$ws.Range('E27').Value = '''  +0.10%  '

# Row 28
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').Value = '''18.57'
$ws.Range('E28').Value = '''  +1.97%  '

# Row 29
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').Value = '''5.329'
$ws.Range('E29').Value = '''  +1.43%  '

# Row 30
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').Value = '''118.20'
$ws.Range('E30').Value = '''  +1.20%  '

# Row 31
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').Value = '''0.08839'
$ws.Range('E31').Value = '''  +1.43%  '

# Row 32
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').Value = '''0.7688'
$ws.Range('E32').Value = '''  +4.25%  '

# Row 33
$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').Value = '''1.180'
$ws.Range('E33').Value = '''  -0.07%  '

# Row 34
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').Value = '''4.563'
$ws.Range('E34').Value = '''  +2.97%  '

# Row 35
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').Value = '''2.888'
$ws.Range('E35').Value = '''  -0.46%  '

# Row 36
$ws.Range('B36').Value = 'Frax'
$ws.Range('C36').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D36').Value = '''0.9999'
$ws.Range('E36').Value = '''  -0.10%  '

# Row 37
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').Value = '''1.103'
$ws.Range('E37').Value = '''  +0.89%  '

# Row 38
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '''0.01985'
$ws.Range('E38').Value = '''  +2.10%  '

# Row 39
$ws.Range('D39').Value = '''0.05282'
$ws.Range('E39').Value = '''  +0.76%  '

# Row 40
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = '''7.436'
$ws.Range('E40').Value = '''  +1.34%  '

# Row 41
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value = '''0.5322'
$ws.Range('E41').Value = '''  +4.91%  '

# Row 42
$ws.Range('D42').Value = '''0.1718'
$ws.Range('E42').Value = '''  +1.75%  '

# Row 43
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').Value = '''2.845'
$ws.Range('E43').Value = '''  -0.95%  '

# Row 44
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').Value = '''2.195'
$ws.Range('E44').Value = '''  +10.50%  '

# Row 45
$ws.Range('B45').Value = 'Aptos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D45').Value = '''8.724'
$ws.Range('E45').Value = '''  +1.85%  '

# Row 46
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = '''0.5057'
$ws.Range('E46').Value = '''  +6.68%  '

# Row 47
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').Value = '''10.58'
$ws.Range('E47').Value = '''  +0.83%  '

# Row 48
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = '''1.707'
$ws.Range('E48').Value = '''  +3.40%  '

# Row 49
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').Value = '''105.68'
$ws.Range('E49').Value = '''  -0.09%  '

# Row 50
$ws.Range('B50').Value = 'PaxDollar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D50').Value = '''0.9996'
$ws.Range('E50').Value = '''  -0.10%  '

# Row 51
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '''0.06371'
$ws.Range('E51').Value = '''  +0.77%  '
